$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 93: LEM-218-00-27KS / LED-218-S00-27
$ws.Range("A93").Value = "LEM-218-00-27KS"
$ws.Range("B93").Value = "LED-218-S00-27"
$ws.Range("C93").Value = 42.3491
$ws.Range("D93").Value = 39.85

# Row 94: LEM-293-00-27KH / LED-293-H00-27 (A94 centered like many existing rows)
$ws.Range("A94").Value = "LEM-293-00-27KH"
$ws.Range("A94").HorizontalAlignment = -4108
$ws.Range("A94").VerticalAlignment = -4108
$ws.Range("B94").Value = "LED-293-H00-27"
$ws.Range("C94").Value = 25.859000000000002
$ws.Range("D94").Value = 23.95

# Row 95: LEM-326-00-30KH / LED-326-H00-30
$ws.Range("A95").Value = "LEM-326-00-30KH"
$ws.Range("B95").Value = "LED-326-H00-30"
$ws.Range("C95").Value = 11.808999999999999
$ws.Range("D95").Value = 9.9

# B93:B95 get centered alignment
$ws.Range("B93:B95").HorizontalAlignment = -4108
$ws.Range("B93:B95").VerticalAlignment = -4108

# Row 96: a blank row with currency number format in A96
$ws.Range("A96").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# Update the sheet view to reflect the new selection/scroll position
$ws.Range("D95").Select()
$ws.Application.ActiveWindow.ScrollRow = 82
